$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: S1_TC_T4 - Get watchlists (repeat of T2 style test)
$ws.Range("A5").Value = "S1_TC_T4"
$ws.Range("B5").Value = "Get watchlists"
$ws.Range("C5").Value = "1PCITATIONS"
$ws.Range("D5").Value = "/lists/watchlist"
$ws.Range("E5").Value = "GET"
$ws.Range("F5").Value = "x-1p-user=b58af128-88d7-4a62-85d0-0ff28f49a9c3"
$ws.Range("L5").Value = "PASS"

# Row 6: S1_TC_T5 - Delete item from watchlist
$ws.Range("A6").Value = "S1_TC_T5"
$ws.Range("B6").Value = "Delete item from watchlist"
$ws.Range("C6").Value = "1PCITATIONS"
$ws.Range("D6").Value = "/lists/watchlist/(S1_TC_T1_hits.hits._id)"
$ws.Range("E6").Value = "DELETE"
$ws.Range("F6").Value = "x-1p-user=b58af128-88d7-4a62-85d0-0ff28f49a9c3"
$ws.Range("I6").Value = "S1_TC_T3"
$ws.Range("J6").Value = "status=200"
$ws.Range("L6").Value = "PASS"

# Row 7: S1_TC_T6 - Get watchlists (repeat again)
$ws.Range("A7").Value = "S1_TC_T6"
$ws.Range("B7").Value = "Get watchlists"
$ws.Range("C7").Value = "1PCITATIONS"
$ws.Range("D7").Value = "/lists/watchlist"
$ws.Range("E7").Value = "GET"
$ws.Range("F7").Value = "x-1p-user=b58af128-88d7-4a62-85d0-0ff28f49a9c3"
$ws.Range("L7").Value = "PASS"

# Update the selected/active cell to match the authored state
$ws.Range("F7").Select() | Out-Null
